# Append: 2025-11-10 06:38 JST
# Update the "取得日時" (acquisition timestamp) column for the existing
# case rows (A2:A8) on the "ランサーズ" sheet to reflect the latest fetch time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-10 06:38:04"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
